$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 203.66667
$ws.Cells.Item(41, 10).Value = 203.66667
$ws.Cells.Item(41, 12).Value = 203.66667
$ws.Cells.Item(41, 14).Value = -1083.66667
$ws.Cells.Item(49, 8).Value = 400
$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 13).ClearContents()
$ws.Cells.Item(74, 8).Value = 4146.75
$ws.Cells.Item(74, 9).Value = 3480
$ws.Cells.Item(74, 11).Value = 3480
$ws.Cells.Item(74, 13).Value = -2544
$ws.Cells.Item(77, 8).Value = 4146.75
$ws.Cells.Item(77, 9).Value = 3480
$ws.Cells.Item(77, 11).Value = 17400
$ws.Cells.Item(77, 13).Value = -12720
$ws.Cells.Item(86, 8).Value = 228601400
$ws.Cells.Item(86, 9).Value = 266700980
$ws.Cells.Item(86, 11).Value = 266700980
$ws.Cells.Item(86, 13).Value = -266699857
$ws.Cells.Item(89, 8).Value = 228601400
$ws.Cells.Item(89, 9).Value = 266700980
$ws.Cells.Item(89, 11).Value = 1333504900
$ws.Cells.Item(89, 13).Value = -1333499284
$ws.Cells.Item(116, 8).Value = 2822.5
$ws.Cells.Item(116, 9).Value = 2521
$ws.Cells.Item(116, 11).Value = 2521
$ws.Cells.Item(116, 13).Value = 921
$ws.Cells.Item(132, 8).Value = 5282.1724
$ws.Cells.Item(132, 9).Value = 4062.9546
$ws.Cells.Item(132, 10).Value = 9114
$ws.Cells.Item(132, 11).Value = 12188.8638
$ws.Cells.Item(132, 12).Value = 27342
$ws.Cells.Item(132, 13).Value = -9658.863799999999
$ws.Cells.Item(132, 14).Value = -32402
$ws.Cells.Item(137, 8).Value = 2737.1396
$ws.Cells.Item(137, 10).Value = 5346.6665
$ws.Cells.Item(137, 12).Value = 16039.9995
$ws.Cells.Item(137, 14).Value = -21139.9995
$ws.Cells.Item(138, 8).Value = 2816.3333
$ws.Cells.Item(138, 9).Value = 2946.4
$ws.Cells.Item(138, 10).Value = 2686.2666
$ws.Cells.Item(138, 11).Value = 8839.200000000001
$ws.Cells.Item(138, 12).Value = 8058.7998
$ws.Cells.Item(138, 13).Value = -3699.200000000001
$ws.Cells.Item(138, 14).Value = -18338.7998
$ws.Cells.Item(139, 8).Value = 43583.332
$ws.Cells.Item(139, 10).Value = 45375
$ws.Cells.Item(139, 12).Value = 45375
$ws.Cells.Item(139, 14).Value = -55655
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 538342.3
$ws.Cells.Item(32, 9).Value = 656319.5600000001
$ws.Cells.Item(32, 10).Value = 17276.166
$ws.Cells.Item(32, 11).Value = 656319.5600000001
$ws.Cells.Item(32, 12).Value = 17276.166
$ws.Cells.Item(32, 13).Value = -656032.5600000001
$ws.Cells.Item(32, 14).Value = -17850.166
$ws.Cells.Item(61, 8).Value = 4909.9
$ws.Cells.Item(61, 9).Value = 5766.6665
$ws.Cells.Item(61, 10).Value = 4542.7144
$ws.Cells.Item(61, 11).Value = 5766.6665
$ws.Cells.Item(61, 12).Value = 4542.7144
$ws.Cells.Item(61, 13).Value = -5554.6665
$ws.Cells.Item(61, 14).Value = -4966.7144
$ws.Cells.Item(74, 8).Value = 1300.6666
$ws.Cells.Item(74, 9).Value = 850.8125
$ws.Cells.Item(74, 11).Value = 850.8125
$ws.Cells.Item(74, 13).Value = 23.1875
$ws.Cells.Item(77, 8).Value = 1300.6666
$ws.Cells.Item(77, 9).Value = 850.8125
$ws.Cells.Item(77, 11).Value = 4254.0625
$ws.Cells.Item(77, 13).Value = 113.9375
$ws.Cells.Item(132, 8).Value = 3276.672
$ws.Cells.Item(132, 9).Value = 2330.0222
$ws.Cells.Item(132, 11).Value = 6990.0666
$ws.Cells.Item(132, 13).Value = -4460.0666
$ws.Cells.Item(136, 8).Value = 4909.9
$ws.Cells.Item(136, 9).Value = 5766.6665
$ws.Cells.Item(136, 10).Value = 4542.7144
$ws.Cells.Item(136, 11).Value = 17299.9995
$ws.Cells.Item(136, 12).Value = 13628.1432
$ws.Cells.Item(136, 13).Value = -14749.9995
$ws.Cells.Item(136, 14).Value = -18728.1432
$ws.Cells.Item(139, 8).Value = 66927.5
$ws.Cells.Item(139, 9).Value = 90000
$ws.Cells.Item(139, 10).Value = 59236.668
$ws.Cells.Item(139, 11).Value = 90000
$ws.Cells.Item(139, 12).Value = 59236.668
$ws.Cells.Item(139, 13).Value = -84860
$ws.Cells.Item(139, 14).Value = -69516.66800000001
$ws.Cells.Item(140, 8).Value = 107499.75
$ws.Cells.Item(140, 10).Value = 107499.75
$ws.Cells.Item(140, 12).Value = 107499.75
$ws.Cells.Item(140, 14).Value = -117859.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3604.96
$ws.Cells.Item(134, 9).Value = 3407.9285
$ws.Cells.Item(134, 10).Value = 3855.7273
$ws.Cells.Item(134, 11).Value = 10223.7855
$ws.Cells.Item(134, 12).Value = 11567.1819
$ws.Cells.Item(134, 13).Value = -7688.7855
$ws.Cells.Item(134, 14).Value = -16637.1819
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6830.3125
$ws.Cells.Item(31, 9).Value = 1267.88
$ws.Cells.Item(31, 10).Value = 12876.435
$ws.Cells.Item(31, 11).Value = 1267.88
$ws.Cells.Item(31, 12).Value = 12876.435
$ws.Cells.Item(31, 13).Value = -972.8800000000001
$ws.Cells.Item(31, 14).Value = -13466.435
$ws.Cells.Item(34, 8).Value = 6830.3125
$ws.Cells.Item(34, 9).Value = 1267.88
$ws.Cells.Item(34, 10).Value = 12876.435
$ws.Cells.Item(34, 11).Value = 1267.88
$ws.Cells.Item(34, 12).Value = 12876.435
$ws.Cells.Item(34, 13).Value = -1065.88
$ws.Cells.Item(34, 14).Value = -13280.435
$ws.Cells.Item(58, 8).Value = 1670.9375
$ws.Cells.Item(58, 9).Value = 1386.091
$ws.Cells.Item(58, 11).Value = 1386.091
$ws.Cells.Item(58, 13).Value = -1183.091
$ws.Cells.Item(132, 8).Value = 6412267
$ws.Cells.Item(132, 9).Value = 1792.7222
$ws.Cells.Item(132, 11).Value = 5378.1666
$ws.Cells.Item(132, 13).Value = -2848.1666
$ws.Cells.Item(134, 8).Value = 3923.4614
$ws.Cells.Item(134, 9).Value = 2714.2856
$ws.Cells.Item(134, 10).Value = 5334.1665
$ws.Cells.Item(134, 11).Value = 8142.8568
$ws.Cells.Item(134, 12).Value = 16002.4995
$ws.Cells.Item(134, 13).Value = -5607.8568
$ws.Cells.Item(134, 14).Value = -21072.4995
$ws.Cells.Item(136, 8).Value = 1670.9375
$ws.Cells.Item(136, 9).Value = 1386.091
$ws.Cells.Item(136, 11).Value = 4158.272999999999
$ws.Cells.Item(136, 13).Value = -1608.272999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 685.5294
$ws.Cells.Item(5, 9).Value = 678.375
$ws.Cells.Item(5, 11).Value = 2035.125
$ws.Cells.Item(5, 13).Value = -1923.125
$ws.Cells.Item(20, 8).Value = 998.3333
$ws.Cells.Item(20, 9).Value = 993.3333
$ws.Cells.Item(20, 11).Value = 2979.9999
$ws.Cells.Item(20, 13).Value = -2752.9999
$ws.Cells.Item(22, 8).Value = 1186.2162
$ws.Cells.Item(22, 9).Value = 961.25
$ws.Cells.Item(22, 10).Value = 1248.2759
$ws.Cells.Item(22, 11).Value = 2883.75
$ws.Cells.Item(22, 12).Value = 3744.8277
$ws.Cells.Item(22, 13).Value = -2714.75
$ws.Cells.Item(22, 14).Value = -4082.8277
$ws.Cells.Item(27, 8).Value = 1186.2162
$ws.Cells.Item(27, 9).Value = 961.25
$ws.Cells.Item(27, 10).Value = 1248.2759
$ws.Cells.Item(27, 11).Value = 2883.75
$ws.Cells.Item(27, 12).Value = 3744.8277
$ws.Cells.Item(27, 13).Value = -2781.75
$ws.Cells.Item(27, 14).Value = -3948.8277
$ws.Cells.Item(32, 8).Value = 9262343
$ws.Cells.Item(32, 10).Value = 9526930
$ws.Cells.Item(32, 12).Value = 28580790
$ws.Cells.Item(32, 14).Value = -28581356
$ws.Cells.Item(110, 8).Value = 13350.583
$ws.Cells.Item(110, 10).Value = 13418.044
$ws.Cells.Item(110, 12).Value = 40254.132
$ws.Cells.Item(110, 14).Value = -48434.132
$ws.Cells.Item(135, 8).Value = 685.5294
$ws.Cells.Item(135, 9).Value = 678.375
$ws.Cells.Item(135, 11).Value = 6105.375
$ws.Cells.Item(135, 13).Value = -3570.375
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3035.7827
$ws.Cells.Item(132, 9).Value = 3028.182
$ws.Cells.Item(132, 11).Value = 9084.545999999998
$ws.Cells.Item(132, 13).Value = -6554.545999999998
$ws.Cells.Item(141, 8).Value = 46639.5
$ws.Cells.Item(141, 10).Value = 46639.5
$ws.Cells.Item(141, 12).Value = 46639.5
$ws.Cells.Item(141, 14).Value = -56999.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 900
$ws.Cells.Item(46, 10).Value = 1033.3334
$ws.Cells.Item(46, 12).Value = 1033.3334
$ws.Cells.Item(46, 14).Value = -1409.3334
$ws.Cells.Item(128, 8).Value = 37966.668
$ws.Cells.Item(128, 10).Value = 37966.668
$ws.Cells.Item(128, 12).Value = 37966.668
$ws.Cells.Item(128, 14).Value = -47926.668
$ws.Cells.Item(132, 8).Value = 2538.5
$ws.Cells.Item(132, 9).Value = 1927.9412
$ws.Cells.Item(132, 11).Value = 5783.8236
$ws.Cells.Item(132, 13).Value = -3253.8236
$ws.Cells.Item(138, 8).Value = 38286
$ws.Cells.Item(138, 10).Value = 38286
$ws.Cells.Item(138, 12).Value = 38286
$ws.Cells.Item(138, 14).Value = -48566
$ws.Cells.Item(141, 8).Value = 80215
$ws.Cells.Item(141, 10).Value = 80215
$ws.Cells.Item(141, 12).Value = 80215
$ws.Cells.Item(141, 14).Value = -90575
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 5378033
$ws.Cells.Item(132, 9).Value = 1272.625
$ws.Cells.Item(132, 10).Value = 23812640
$ws.Cells.Item(132, 11).Value = 3817.875
$ws.Cells.Item(132, 12).Value = 71437920
$ws.Cells.Item(132, 13).Value = -1287.875
$ws.Cells.Item(132, 14).Value = -71442980
$ws.Cells.Item(136, 8).Value = 2332.2126
$ws.Cells.Item(136, 9).Value = 2121.3
$ws.Cells.Item(136, 11).Value = 6363.900000000001
$ws.Cells.Item(136, 13).Value = -3813.900000000001
